$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate()

# fix_demand_to_representative_year: FALSE -> TRUE
$ws.Range("B28").Value2 = $true

# fix_profiles_to_representative_year: FALSE -> TRUE
$ws.Range("B29").Value2 = $true

# Representative year: 2004 -> 2015
$ws.Range("B30").Value2 = 2015

# capacity_remuneration_mechanism: strategic_reserve_ger -> none
$ws.Range("B44").Value2 = "none"

# scenarioWeatheryearsExcel file name: 40weatherYears2050TNO.xlsx -> 40weatherYears2050TNO-S1.xlsx
$ws.Range("B32").Value2 = "40weatherYears2050TNO-S1.xlsx"

# Recalculate dependent formulas in column C
$wb.Application.CalculateFull()

# Update the visible selection/scroll position to match where the author ended up
$ws.Range("A30").Select()
$excel.ActiveWindow.ScrollRow = 30
$ws.Range("C51").Select()
